$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 14 ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A14").Value = "Kun jij 5 pakken A4-papier bestellen?"
$ws.Range("B14").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$ws.Range("C14").Value = "Hoi Johan`nKun jij 5 pakken A4-papier bestellen?`nMarc`nSent using {0}"
$ws.Range("D14").Value = "Bestelling / Levering"
$ws.Range("E14").Value = "Beste Marc,`nBedankt voor je e-mail. Ik zorg ervoor dat er 5 pakken A4-papier worden besteld. Heb je nog specifieke voorkeuren voor het merk of type papier? Laat het me weten als je nog andere vragen hebt.`nMet vriendelijke groet,`nJohan"
$ws.Range("F14").Value = "2025-06-26 21:03:25"
$ws.Range("G14").Value = "Ja"
$ws.Range("H14").Value = "Nee"
$ws.Range("I14").Value = "Ja"

# Entering the multi-line text above makes the engine auto-expand the row
# height; AutoFit restores it to the sheet default (matches the other rows,
# none of which carry an explicit height).
$ws.Rows.Item(14).AutoFit()

# --- Extend the conditional formatting ranges to include the new row ---
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = "$col" + "2:$col" + "13"
    $newRange = "$col" + "2:$col" + "14"
    $fcs = $ws.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($newRange))
    }
}

# --- Sheet "Dashboard": bump the "Bestelling / Levering" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 9
